$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.795.21"
$ws.Range("E2").Value = "  +4.03%  "
$ws.Range("D3").Value = "2.275.88"
$ws.Range("E3").Value = "  +4.60%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'251.30"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "'0.627"
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").Value = "'71.58"
$ws.Range("E7").Value = "  +8.09%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").Value = "'0.651"
$ws.Range("E9").Value = "  +13.25%  "
$ws.Range("D10").Value = "'38.45"
$ws.Range("E10").Value = "  +5.30%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.0967"
$ws.Range("E11").Value = "  +3.37%  "
$ws.Range("B12").Value = "OKB"
$ws.Range("C12").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D12").Value = "'59.69"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("E13").Value = "  +6.84%  "
$ws.Range("D14").Value = "'0.106"
$ws.Range("E14").Value = "  +2.49%  "
$ws.Range("D15").Value = "2.618.38"
$ws.Range("E15").Value = "  +4.66%  "
$ws.Range("D16").Value = "'14.94"
$ws.Range("E16").Value = "  +4.62%  "
$ws.Range("D17").Value = "'0.883"
$ws.Range("E17").Value = "  +3.53%  "
$ws.Range("D18").Value = "2.268.30"
$ws.Range("E18").Value = "  +3.95%  "
$ws.Range("D19").Value = "42.766.21"
$ws.Range("E19").Value = "  +4.03%  "
$ws.Range("E20").Value = "  +7.03%  "
$ws.Range("D21").Value = "'6.32"
$ws.Range("E21").Value = "  +3.91%  "
$ws.Range("D22").Value = "'73.40"
$ws.Range("E22").Value = "  +2.29%  "
$ws.Range("D23").Value = "'234.44"
$ws.Range("E23").Value = "  +2.36%  "
$ws.Range("D24").Value = "'2.13"
$ws.Range("E24").Value = "  +3.80%  "
$ws.Range("D25").Value = "'3.91"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").Value = "'11.62"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -1.22%  "
$ws.Range("D30").Value = "'2.12"
$ws.Range("E30").Value = "  +4.44%  "
$ws.Range("D31").Value = "'168.39"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("D32").Value = "'21.37"
$ws.Range("E32").Value = "  +5.56%  "
$ws.Range("D33").Value = "'6.40"
$ws.Range("E33").Value = "  +10.58%  "
$ws.Range("E34").Value = "  +6.76%  "
$ws.Range("D35").Value = "'0.0800"
$ws.Range("E35").Value = "  +4.78%  "
$ws.Range("D36").Value = "'30.68"
$ws.Range("E36").Value = "  +25.00%  "
$ws.Range("E37").Value = "  +4.13%  "
$ws.Range("D38").Value = "'4.65"
$ws.Range("E38").Value = "  +16.58%  "
$ws.Range("E39").Value = "  +4.86%  "
$ws.Range("D40").Value = "'0.0311"
$ws.Range("E40").Value = "  +1.00%  "
$ws.Range("D41").Value = "'2.32"
$ws.Range("E41").Value = "  +4.86%  "
$ws.Range("D42").Value = "'13.21"
$ws.Range("E42").Value = "  +16.22%  "
$ws.Range("D43").Value = "'5.87"
$ws.Range("E43").Value = "  +6.76%  "
$ws.Range("E44").Value = "  +11.53%  "
$ws.Range("B45").Value = "FTXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D45").Value = "'5.00"
$ws.Range("E45").Value = "  -5.56%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.17"
$ws.Range("E46").Value = "  +7.68%  "
$ws.Range("D47").Value = "'61.43"
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").Value = "'0.102"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E49").Value = "  +3.76%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("E51").Value = "  +4.03%  "
